$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$c = $ws.Range("D2")
$c.NumberFormat = "@"
$c.Value = "30.501.81"
$c.Style = "Normal"
$ws.Range("E2").Value = "  +0.20%  "
$c = $ws.Range("D3")
$c.NumberFormat = "@"
$c.Value = "1.914.70"
$c.Style = "Normal"
$ws.Range("E3").Value = "  -0.15%  "
$c = $ws.Range("D4")
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Range("E4").Value = "  -0.03%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "245.22"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +0.66%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "0.9989"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -0.06%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.4803"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +2.20%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.2885"
$c.Style = "Normal"
$ws.Range("E8").Value = "  +0.85%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.06734"
$c.Style = "Normal"
$ws.Range("E9").Value = "  -1.48%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "110.87"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +0.45%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "19.21"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.08%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "1.913.05"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +1.27%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "0.07551"
$c.Style = "Normal"
$ws.Range("E13").Value = "  -2.41%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "5.249"
$c.Style = "Normal"
$ws.Range("E14").Value = "  -1.04%  "
$c = $ws.Range("D15")
$c.NumberFormat = "@"
$c.Value = "0.6682"
$c.Style = "Normal"
$ws.Range("E15").Value = "  +1.39%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "300.84"
$c.Style = "Normal"
$ws.Range("E16").Value = "  +1.66%  "
$c = $ws.Range("D17")
$c.NumberFormat = "@"
$c.Value = "30.501.28"
$c.Style = "Normal"
$ws.Range("E17").Value = "  +0.22%  "
$c = $ws.Range("D18")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E18").Value = "  -0.02%  "
$ws.Range("E19").Value = "  +0.21%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "0.000007575"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -0.80%  "
$ws.Range("B21").Value = "Uniswap"
$ws.Range("C21").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "5.501"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +4.64%  "
$ws.Range("B22").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C22").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "2.162.91"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +1.04%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E23").Value = "  -0.18%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "6.403"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +2.82%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "9.481"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +1.23%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "164.27"
$c.Style = "Normal"
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "20.60"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -6.04%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "2.099"
$c.Style = "Normal"
$ws.Range("E28").Value = "  +0.48%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.1072"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +0.21%  "
$ws.Range("E30").Value = "  +2.64%  "
$ws.Range("E31").Value = "  -0.28%  "
$ws.Range("E32").Value = "  +0.87%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.04969"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -1.59%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "0.7303"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -0.90%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.136"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -1.78%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "0.02053"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -0.67%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "0.9994"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +0.10%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "2.727"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -0.56%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.670"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -0.28%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "110.93"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +1.19%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "2.011"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -2.38%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.4410"
$c.Style = "Normal"
$ws.Range("E42").Value = "  +3.45%  "
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.8637"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -1.06%  "
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "5.882"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +0.81%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "0.9988"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -0.03%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "68.47"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +1.36%  "
$ws.Range("B47").Value = "BitcoinSV"
$ws.Range("C47").Value = "https://coinranking.com/coin/VcMY11NONHSA0+bitcoinsv-bsv"
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "50.04"
$c.Style = "Normal"
$ws.Range("E47").Value = "  -2.55%  "
$ws.Range("B48").Value = "Aptos"
$ws.Range("C48").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "7.302"
$c.Style = "Normal"
$ws.Range("E48").Value = "  +1.40%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "9.254"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +0.38%  "
$ws.Range("E50").Value = "  +1.53%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "0.2545"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +3.28%  "
